$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Warehouse")

# Row 29: DimAgent.WorkingDuration -> DimAgent.HireDate, destination type int -> datetime
$ws.Range("A29").Value = "DWRedwood.dbo.DimAgent.HireDate"
$ws.Range("E29").Value = "datetime"

# Row 30: DimAgent.Age -> DimAgent.BirthDate, destination type int -> datetime
$ws.Range("A30").Value = "DWRedwood.dbo.DimAgent.BirthDate"
$ws.Range("E30").Value = "datetime"
